$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 529.36365
$ws.Range("J17").Value = 529.36365
$ws.Range("L17").Value = 1588.09095
$ws.Range("N17").Value = -1924.09095
$ws.Range("H87").Value = 35333.332
$ws.Range("J87").Value = 35333.332
$ws.Range("L87").Value = 35333.332
$ws.Range("N87").Value = -37829.332
$ws.Range("H90").Value = 35333.332
$ws.Range("J90").Value = 35333.332
$ws.Range("L90").Value = 105999.996
$ws.Range("N90").Value = -118479.996
$ws.Range("H97").Value = 2198.3333
$ws.Range("J97").Value = 2198.3333
$ws.Range("L97").Value = 6594.999899999999
$ws.Range("N97").Value = -7586.999899999999
$ws.Range("H107").Value = 9804378
$ws.Range("I107").Value = 16666933
$ws.Range("J107").Value = 728.4286
$ws.Range("K107").Value = 16666933
$ws.Range("L107").Value = 728.4286
$ws.Range("M107").Value = -16665013
$ws.Range("N107").Value = -4568.4286
$ws.Range("H112").Value = 2711792.8
$ws.Range("J112").Value = 2925811.8
$ws.Range("L112").Value = 8777435.399999999
$ws.Range("N112").Value = -8779651.399999999
$ws.Range("H138").Value = 2851262
$ws.Range("I138").Value = 209803.5
$ws.Range("J138").Value = 11907691
$ws.Range("K138").Value = 629410.5
$ws.Range("L138").Value = 35723073
$ws.Range("M138").Value = -624270.5
$ws.Range("N138").Value = -35733353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1073.6
$ws.Range("I32").Value = 984.022
$ws.Range("J32").Value = 1979.3334
$ws.Range("K32").Value = 984.022
$ws.Range("L32").Value = 1979.3334
$ws.Range("M32").Value = -697.022
$ws.Range("N32").Value = -2553.3334
$ws.Range("H45").Value = 1522
$ws.Range("J45").Value = 1729.8334
$ws.Range("L45").Value = 1729.8334
$ws.Range("N45").Value = -2483.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 5611.1113
$ws.Range("J19").Value = 5611.1113
$ws.Range("L19").Value = 5611.1113
$ws.Range("N19").Value = -5957.1113
$ws.Range("H107").Value = 2006.1578
$ws.Range("I107").Value = 2460.2
$ws.Range("J107").Value = 1501.6666
$ws.Range("K107").Value = 2460.2
$ws.Range("L107").Value = 1501.6666
$ws.Range("M107").Value = -540.1999999999998
$ws.Range("N107").Value = -5341.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 797.2308
$ws.Range("I5").Value = 395
$ws.Range("J5").Value = 870.36365
$ws.Range("K5").Value = 395
$ws.Range("L5").Value = 870.36365
$ws.Range("M5").Value = -283
$ws.Range("N5").Value = -1094.36365
$ws.Range("H12").Value = 10276.667
$ws.Range("J12").Value = 14950
$ws.Range("L12").Value = 14950
$ws.Range("N12").Value = -15290
$ws.Range("H31").Value = 5273.7617
$ws.Range("I31").Value = 5262.5
$ws.Range("J31").Value = 5276.4116
$ws.Range("K31").Value = 5262.5
$ws.Range("L31").Value = 5276.4116
$ws.Range("M31").Value = -4967.5
$ws.Range("N31").Value = -5866.4116
$ws.Range("H34").Value = 5273.7617
$ws.Range("I34").Value = 5262.5
$ws.Range("J34").Value = 5276.4116
$ws.Range("K34").Value = 5262.5
$ws.Range("L34").Value = 5276.4116
$ws.Range("M34").Value = -5060.5
$ws.Range("N34").Value = -5680.4116
$ws.Range("H62").Value = 2987.5
$ws.Range("I62").Value = 2975
$ws.Range("K62").Value = 2975
$ws.Range("M62").Value = -2351
$ws.Range("H65").Value = 2987.5
$ws.Range("I65").Value = 2975
$ws.Range("K65").Value = 14875
$ws.Range("M65").Value = -11755
$ws.Range("H99").Value = 7764.6
$ws.Range("I99").Value = 9421.5
$ws.Range("J99").Value = 6660
$ws.Range("K99").Value = 9421.5
$ws.Range("L99").Value = 6660
$ws.Range("M99").Value = -7923.5
$ws.Range("N99").Value = -9656
$ws.Range("H126").Value = 7764.6
$ws.Range("I126").Value = 9421.5
$ws.Range("J126").Value = 6660
$ws.Range("K126").Value = 28264.5
$ws.Range("L126").Value = 19980
$ws.Range("M126").Value = -25794.5
$ws.Range("N126").Value = -24920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15873950
$ws.Range("I131").Value = 83333704
$ws.Range("J131").Value = 1066.804
$ws.Range("K131").Value = 250001112
$ws.Range("L131").Value = 3200.412
$ws.Range("M131").Value = -249996072
$ws.Range("N131").Value = -13280.412
$ws.Range("H138").Value = 3151.5
$ws.Range("I138").Value = 2484.5454
$ws.Range("J138").Value = 3966.6667
$ws.Range("K138").Value = 7453.6362
$ws.Range("L138").Value = 11900.0001
$ws.Range("M138").Value = -2313.6362
$ws.Range("N138").Value = -22180.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H102").Value = 1498.6666
$ws.Range("I102").Value = 975
$ws.Range("K102").Value = 975
$ws.Range("M102").Value = 647
$ws.Range("H108").Value = 24884
$ws.Range("J108").Value = 24884
$ws.Range("L108").Value = 24884
$ws.Range("N108").Value = -32564
$ws.Range("H110").Value = 49602
$ws.Range("J110").Value = 49602
$ws.Range("L110").Value = 49602
$ws.Range("N110").Value = -57782
$ws.Range("H126").Value = 2841.1765
$ws.Range("I126").Value = 5266.6665
$ws.Range("J126").Value = 2321.4285
$ws.Range("K126").Value = 15799.9995
$ws.Range("L126").Value = 6964.2855
$ws.Range("M126").Value = -13329.9995
$ws.Range("N126").Value = -11904.2855
$ws.Range("H141").Value = 30834.875
$ws.Range("J141").Value = 30834.875
$ws.Range("L141").Value = 30834.875
$ws.Range("N141").Value = -41194.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2671.4285
$ws.Range("I7").Value = 2233.3333
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2233.3333
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2121.3333
$ws.Range("N7").Value = -3224
$ws.Range("H40").Value = 2519.9333
$ws.Range("I40").Value = 2461.8462
$ws.Range("J40").Value = 2897.5
$ws.Range("K40").Value = 2461.8462
$ws.Range("L40").Value = 2897.5
$ws.Range("M40").Value = -2325.8462
$ws.Range("N40").Value = -3169.5
$ws.Range("H69").Value = 40163
$ws.Range("J69").Value = 40163
$ws.Range("L69").Value = 40163
$ws.Range("N69").Value = -41785
$ws.Range("H72").Value = 40163
$ws.Range("J72").Value = 40163
$ws.Range("L72").Value = 120489
$ws.Range("N72").Value = -128601
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws.Range("H126").Value = 2671.4285
$ws.Range("I126").Value = 2233.3333
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6699.999899999999
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4229.999899999999
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1150.7693
$ws.Range("I126").Value = 969.6
$ws.Range("J126").Value = 1264
$ws.Range("K126").Value = 2908.8
$ws.Range("L126").Value = 3792
$ws.Range("M126").Value = -438.8000000000002
$ws.Range("N126").Value = -8732
$ws.Range("H132").Value = 65239.02
$ws.Range("I132").Value = 54450.406
$ws.Range("J132").Value = 101528
$ws.Range("K132").Value = 163351.218
$ws.Range("L132").Value = 304584
$ws.Range("M132").Value = -160821.218
$ws.Range("N132").Value = -309644
$ws.Range("H136").Value = 32373.969
$ws.Range("I136").Value = 21885.766
$ws.Range("J136").Value = 61370.766
$ws.Range("K136").Value = 65657.298
$ws.Range("L136").Value = 184112.298
$ws.Range("M136").Value = -63107.298
$ws.Range("N136").Value = -189212.298
